$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# ---------------------------------------------------------------------------
# Table 1 (rows 1-11): update quantities of cervejas/refrigerante/espetinhos
# ---------------------------------------------------------------------------
$ws.Range("B3").Value = 1
$ws.Range("B4").Value = 1
$ws.Range("B5").Value = 1

# Rework the tip ("10% GARCON") row into a "Subtotal sem 10%" row, and the
# formulas that depend on it.
$ws.Range("B10").Formula = "=(B7+B8+B9)"
$ws.Range("C10").Value = "Subtotal sem 10%"

$ws.Range("B11").Formula = "=B10+((B10*10)/100)"

# ---------------------------------------------------------------------------
# Table 2 (rows 14-24)
# ---------------------------------------------------------------------------
$ws.Range("B16").Value = 2
$ws.Range("B17").Value = 2
$ws.Range("B18").Value = 2

$ws.Range("B20").Value = 8
$ws.Range("B21").Value = 17
$ws.Range("B22").Value = 0
$ws.Range("C23").Value = "Subtotal sem 10%"
$ws.Range("B23").Value = 25
$ws.Range("B24").Value = 27.5

# ---------------------------------------------------------------------------
# Table 3 (rows 27-37)
# ---------------------------------------------------------------------------
# Match row 27's header formatting (C27/D27) to the other tables' header row
# (C1/D1 use the same border styling as C14/D14 should, but the target makes
# C27/D27 mirror C1/D1's styles).
$ws.Range("C1:D1").Copy()
$ws.Range("C27:D27").PasteSpecial(-4122)

$ws.Range("B29").Value = 3
$ws.Range("B30").Value = 3
$ws.Range("B31").Value = 3

$ws.Range("B33").Value = 8
$ws.Range("B34").Value = 25.5
$ws.Range("B35").Value = 0
$ws.Range("C36").Value = "Subtotal sem 10%"
$ws.Range("B36").Value = 33.5
$ws.Range("B37").Value = 36.85

# ---------------------------------------------------------------------------
# Column C width + selection
# ---------------------------------------------------------------------------
$ws.Columns("C").ColumnWidth = 16.7109375

$ws.Range("B8").Select()
